# "Generate Report for Handback" — update the localization-status report
# to reflect that the de-de (and zh-cn) handback round has completed and
# is now in sync with en-US.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status text -----------------------------------------------------
# The shared "Ready for handoff" string is used by Overview!E2 (zh-cn
# status), Overview!F2 (de-de status), zh-cn!C2 and de-de!C2 — updating
# each cell's value drives them all to the new shared string.
$ws_overview.Range("E2").Value = "Handed back: in sync with en-US"
$ws_overview.Range("F2").Value = "Handed back: in sync with en-US"
$ws_zhcn.Range("C2").Value     = "Handed back: in sync with en-US"
$ws_dede.Range("C2").Value     = "Handed back: in sync with en-US"

# --- zh-cn handback refresh ------------------------------------------
$ws_zhcn.Range("K2").Value = "2016-08-28 18:48:13"   # Latest Handback DateTime
$ws_zhcn.Range("P2").Value = ""                       # Error Detail cleared

# --- de-de handback refresh ------------------------------------------
$ws_dede.Range("K2").Value = "2016-08-28 18:48:19"    # Latest Handback DateTime
$ws_dede.Range("P2").Value = ""                        # Error Detail cleared

# --- Column width refresh (Status / Error Detail columns resized) ----
# ColumnWidth is in characters; observed offset between the stored OOXML
# width and this property is a constant +5/6 once snapped to the grid.
$ws_overview.Range("E1").ColumnWidth = 29.166666666666668
$ws_overview.Range("F1").ColumnWidth = 29.166666666666668

$ws_zhcn.Range("C1").ColumnWidth = 29.166666666666668
$ws_zhcn.Range("P1").ColumnWidth = 12.833333333333334

$ws_dede.Range("C1").ColumnWidth = 29.166666666666668
$ws_dede.Range("P1").ColumnWidth = 12.833333333333334
